# Fix issue where Somass (GCL, SPR) age comp columns were not proportions.
$wb = $excel.ActiveWorkbook

# --- 1. Update the "fertilized" column definition text on the metadata sheet ---
$meta = $wb.Worksheets.Item("metadata")
$meta.Range("B13").Value = "Binary variable describing whether (1) or not (0) the CU nursary lake was fertilized in each year. Note that fertilization affects abundances of pre-smolts in year + 1 (e.g. fertilizing a lake in 2010 is expected to affect the fry arising from brood year 2009)."

# --- 2. Fix the Somass (GCL / SPR) stock age-composition columns on the "S-R data" sheet ---
$ws = $wb.Worksheets.Item("S-R data")

$ur = $ws.UsedRange
$lastRow = $ur.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $stock = $ws.Cells.Item($r, 2).Value2

    if ($stock -eq "GCL" -or $stock -eq "SPR") {
        $E = $ws.Cells.Item($r, 5).Value2
        $F = $ws.Cells.Item($r, 6).Value2
        $G = $ws.Cells.Item($r, 7).Value2
        $H = $ws.Cells.Item($r, 8).Value2
        $I = $ws.Cells.Item($r, 9).Value2
        $J = $ws.Cells.Item($r, 10).Value2
        $K = $ws.Cells.Item($r, 11).Value2

        $ageSum = $G + $H + $I + $J

        $ws.Cells.Item($r, 5).Value = $E - $K
        $ws.Cells.Item($r, 6).Value = $F - $K
        $ws.Cells.Item($r, 7).Value = $G / $ageSum
        $ws.Cells.Item($r, 8).Value = $H / $ageSum
        $ws.Cells.Item($r, 9).Value = $I / $ageSum
        $ws.Cells.Item($r, 10).Value = $J / $ageSum
    }
}
